# Appointments.xlsx re-upload:
#  - Appointment time text tidied to upper-case "PM"
#  - second row's appointment re-scheduled to a new date/time
#  - selection cursor left on D2
#  - workbook re-saved from a newer Excel build, which swapped the
#    default "Office" theme (Calibri / blue accents) for the newer
#    "Office 2023" theme (Aptos fonts / refreshed accent colours)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edits -----------------------------------------------------
$ws.Range("D2").Value = "18-Nov-2024 2:00:00 PM"
$ws.Range("D3").Value = "20-Nov-2024 5:00:00 PM"

# --- selection --------------------------------------------------------
$ws.Range("D2").Select() | Out-Null

# --- theme refresh (Office -> Office 2023 "Aptos" theme) --------------
$theme = $wb.Theme

$fontScheme = $theme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Aptos Display"
$fontScheme.MinorFont.Latin = "Aptos Narrow"

$colors = $theme.ThemeColorScheme
# index : MsoThemeColorSchemeIndex -> slot
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
$colors.Colors(3).RGB  = 4270094    # dk2      44546A -> 0E2841
$colors.Colors(4).RGB  = 15263976   # lt2      E7E6E6 -> E8E8E8
$colors.Colors(5).RGB  = 8544277    # accent1  4472C4 -> 156082
$colors.Colors(6).RGB  = 3305961    # accent2  ED7D31 -> E97132
$colors.Colors(7).RGB  = 2386713    # accent3  A5A5A5 -> 196B24
$colors.Colors(8).RGB  = 13999631   # accent4  FFC000 -> 0F9ED5
$colors.Colors(9).RGB  = 9644960    # accent5  5B9BD5 -> A02B93
$colors.Colors(10).RGB = 3057486    # accent6  70AD47 -> 4EA72E
$colors.Colors(11).RGB = 8812614    # hlink    0563C1 -> 467886
$colors.Colors(12).RGB = 8216726    # folHlink 954F72 -> 96607D
